$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update KPI target values (column D, and C2)
$ws.Range("C2").Value = 60
$ws.Range("D3").Value = 10
$ws.Range("D4").Value = 20
$ws.Range("D5").Value = 30
$ws.Range("D7").Value = 30
$ws.Range("D8").Value = 0
$ws.Range("D11").Value = 20
$ws.Range("D12").Value = 0
$ws.Range("D15").Value = 30

# Adjust column widths (values widened slightly vs original)
$ws.Columns.Item(1).ColumnWidth = 29.406462585034
$ws.Columns.Item(2).ColumnWidth = 59.9319727891157
$ws.Columns.Item(3).ColumnWidth = 11.1156462585034
$ws.Columns.Item(4).ColumnWidth = 11.1156462585034
$ws.Columns.Item(5).ColumnWidth = 11.1156462585034
$ws.Columns.Item(6).ColumnWidth = 67.8503401360545

# Update the selected cell
$ws.Range("B23").Select()
